# Updated BGR model - 2025-08-17 19:38
# Re-assigns the lcoe_class (cost-class) label/number among a handful of
# resource rows on the "solar" and "wind" sheets. The numeric resource data
# (capacity, capacity factor, LCOE, etc. in columns M:O) stays attached to
# the same physical row; only the process name (C), description (D), the
# mirrored process name (K) and the lcoe_class number (P) move between rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "solar": spv-BGR_16 cost classes 2/3/4 rotate across rows 4-6
# ---------------------------------------------------------------------
$wsSolar = $wb.Worksheets.Item("solar")

$wsSolar.Range("C4").Value = "e_spv-BGR_16_c4"
$wsSolar.Range("D4").Value = "solar resource -- CF class spv-BGR_16 -- cost class 4"
$wsSolar.Range("K4").Value = "e_spv-BGR_16_c4"
$wsSolar.Range("P4").Value = 4

$wsSolar.Range("C5").Value = "e_spv-BGR_16_c2"
$wsSolar.Range("D5").Value = "solar resource -- CF class spv-BGR_16 -- cost class 2"
$wsSolar.Range("K5").Value = "e_spv-BGR_16_c2"
$wsSolar.Range("P5").Value = 2

$wsSolar.Range("C6").Value = "e_spv-BGR_16_c3"
$wsSolar.Range("D6").Value = "solar resource -- CF class spv-BGR_16 -- cost class 3"
$wsSolar.Range("K6").Value = "e_spv-BGR_16_c3"
$wsSolar.Range("P6").Value = 3

# ---------------------------------------------------------------------
# Sheet "wind": several won-BGR_* cost classes rotate across their rows
# ---------------------------------------------------------------------
$wsWind = $wb.Worksheets.Item("wind")

# won-BGR_25: cost classes 1/2/3 rotate across rows 15-17
$wsWind.Range("C15").Value = "e_won-BGR_25_c2"
$wsWind.Range("D15").Value = "wind resource -- CF class won-BGR_25 -- cost class 2"
$wsWind.Range("K15").Value = "e_won-BGR_25_c2"
$wsWind.Range("P15").Value = 2

$wsWind.Range("C16").Value = "e_won-BGR_25_c3"
$wsWind.Range("D16").Value = "wind resource -- CF class won-BGR_25 -- cost class 3"
$wsWind.Range("K16").Value = "e_won-BGR_25_c3"
$wsWind.Range("P16").Value = 3

$wsWind.Range("C17").Value = "e_won-BGR_25_c1"
$wsWind.Range("D17").Value = "wind resource -- CF class won-BGR_25 -- cost class 1"
$wsWind.Range("K17").Value = "e_won-BGR_25_c1"
$wsWind.Range("P17").Value = 1

# won-BGR_24: cost classes 1/3 swap across rows 18-19
$wsWind.Range("C18").Value = "e_won-BGR_24_c1"
$wsWind.Range("D18").Value = "wind resource -- CF class won-BGR_24 -- cost class 1"
$wsWind.Range("K18").Value = "e_won-BGR_24_c1"
$wsWind.Range("P18").Value = 1

$wsWind.Range("C19").Value = "e_won-BGR_24_c3"
$wsWind.Range("D19").Value = "wind resource -- CF class won-BGR_24 -- cost class 3"
$wsWind.Range("K19").Value = "e_won-BGR_24_c3"
$wsWind.Range("P19").Value = 3

# won-BGR_21: cost classes 1/2 swap across rows 27-28
$wsWind.Range("C27").Value = "e_won-BGR_21_c1"
$wsWind.Range("D27").Value = "wind resource -- CF class won-BGR_21 -- cost class 1"
$wsWind.Range("K27").Value = "e_won-BGR_21_c1"
$wsWind.Range("P27").Value = 1

$wsWind.Range("C28").Value = "e_won-BGR_21_c2"
$wsWind.Range("D28").Value = "wind resource -- CF class won-BGR_21 -- cost class 2"
$wsWind.Range("K28").Value = "e_won-BGR_21_c2"
$wsWind.Range("P28").Value = 2

# won-BGR_17: cost classes 1/2 swap across rows 47-48
$wsWind.Range("C47").Value = "e_won-BGR_17_c1"
$wsWind.Range("D47").Value = "wind resource -- CF class won-BGR_17 -- cost class 1"
$wsWind.Range("K47").Value = "e_won-BGR_17_c1"
$wsWind.Range("P47").Value = 1

$wsWind.Range("C48").Value = "e_won-BGR_17_c2"
$wsWind.Range("D48").Value = "wind resource -- CF class won-BGR_17 -- cost class 2"
$wsWind.Range("K48").Value = "e_won-BGR_17_c2"
$wsWind.Range("P48").Value = 2
